$d = $word.ActiveDocument

# Paragraph 1: "A complicated union can be seen ... (Ref-J6DH3Y)."
$d.Paragraphs(1).Range.Find.Execute("Ref-J6DH3Y", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-u314032", 2)

# Paragraph 3: "Due to FL's underlying capability ... (Ref-DJ49F2)."
$d.Paragraphs(3).Range.Find.Execute("Ref-DJ49F2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f282285", 2)

# Paragraph 4: "Due to the rapidly changing ... (Ref-J7X8K2)."
$d.Paragraphs(4).Range.Find.Execute("Ref-J7X8K2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f815469", 2)

# Paragraph 6: "Economic and inequality isolation ..." has three Ref-J7X8K2 citations
$d.Paragraphs(6).Range.Find.Execute("Ref-J7X8K2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-u617530", 2)

# Paragraph 7: "Nobody in the United States ..." has Ref-A1B2C3 and Ref-D4E5F6
$d.Paragraphs(7).Range.Find.Execute("Ref-A1B2C3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s084268", 2)
$d.Paragraphs(7).Range.Find.Execute("Ref-D4E5F6", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s084268", 2)
